# Test cases with the Bioblank and Diagnosis filters
#
# Rewrites the Cypher queries behind the "startup" control sheet:
#  - Row 2 (CasesTab)      : updates the Cases query and points the stat-query
#                            column at the new combined counts query.
#  - Row 3 (SamplesTab)    : updates the Samples query (pathology replace +
#                            ordering) and points at the new stat query.
#  - Row 4 (FilesTab)      : replaces the Files query with the new version
#                            that adds size formatting / sample linkage.
#  - Row 5 (StudyFilesTab) : replaces the Study Files query with the new
#                            version that adds size formatting / ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New / updated Cypher query bodies (verbatim, single-quoted here-strings so
# that backticks and $ signs inside the Cypher text are kept literal).
# ---------------------------------------------------------------------------

$statsQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
MATCH (r:registration)-->(c)
WHERE r.registration_origin = "CSU ACTR"
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$filesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f)-[*]->(samp:sample)
MATCH (f)-[*]->(c:case)<--(demo:demographic)
MATCH (r:registration)-->(c)
WHERE r.registration_origin = "CSU ACTR"
OPTIONAL MATCH (s:study)<--(c)<--(diag:diagnosis)<-[*]-(samp)
WITH
        f, parent, c, demo, diag, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value, 
        10^precision AS factor,
        units[i] as unit
WITH    
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN 
        coalesce(f.file_name, '') AS `File Name`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_type, '') AS `File Type`,
       CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(samp.sample_id, '') AS `Sample ID`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis
Order By f.file_name LIMIT 100
'@

$studyFilesQuery = @'
MATCH (f:file)-->(s:study)
MATCH (s)<--(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
MATCH (r:registration)-->(c)
WHERE r.registration_origin = "CSU ACTR"
WITH DISTINCT f,  s, c, demo, diag
WITH
        f, c, demo, diag, s,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH    
        f, c, demo, diag, s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
        WITH    
        f,  c, demo, diag, s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '') AS `File Name`,
  coalesce(f.file_type, '') AS `File Type`,
  coalesce("study", '') AS `Association`,
  coalesce(f.file_description, '') AS `Description`,
  coalesce(f.file_format, '') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
  coalesce(s.clinical_study_designation,'') AS `Study Code`
  order by 'File Name' asc
  limit 100
'@

$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis), (r:registration)-->(c)
WHERE r.registration_origin = "CSU ACTR"
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight
RETURN  
       coalesce(c.case_id, '') AS `Case ID`,
       coalesce(s.clinical_study_designation, '') AS `Study Code`,
       coalesce(s.clinical_study_type, '') AS  `Study Type`,
       coalesce(demo.breed, '') AS Breed ,
       coalesce(diag.disease_term, '') AS Diagnosis ,
       coalesce(diag.stage_of_disease, '') AS `Stage of Disease`,
      coalesce(CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END, '') AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
       coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
       coalesce(diag.best_response, '') AS `Response to Treatment`,
       coalesce(co.cohort_description, '') AS `Cohort`
Order by c.case_id LIMIT 100
'@

$samplesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis), (r:registration)-->(c)
WHERE r.registration_origin = "CSU ACTR"
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed,
        coalesce(diag.disease_term,'') AS Diagnosis,
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
       coalesce(replace(samp.specific_sample_pathology,"  "," "),'') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
Order by samp.sample_id LIMIT 100
'@

# Here-strings capture a trailing newline before the closing "'@" marker;
# trim it so the stored text matches the query text exactly.
$statsQuery      = $statsQuery.TrimEnd("`r", "`n")
$filesQuery      = $filesQuery.TrimEnd("`r", "`n")
$studyFilesQuery = $studyFilesQuery.TrimEnd("`r", "`n")
$casesQuery      = $casesQuery.TrimEnd("`r", "`n")
$samplesQuery    = $samplesQuery.TrimEnd("`r", "`n")

# ---------------------------------------------------------------------------
# Apply the updated queries to the worksheet (column B = query, column C =
# the "stat" query that now is shared/identical across all four tabs).
# ---------------------------------------------------------------------------

$ws.Range("B2").Value = $casesQuery        # CasesTab
$ws.Range("C2").Value = $statsQuery

$ws.Range("B3").Value = $samplesQuery      # SamplesTab
$ws.Range("C3").Value = $statsQuery

$ws.Range("B4").Value = $filesQuery        # FilesTab
$ws.Range("C4").Value = $statsQuery

$ws.Range("B5").Value = $studyFilesQuery   # StudyFilesTab
$ws.Range("C5").Value = $statsQuery

# ---------------------------------------------------------------------------
# The longer query text now needs more vertical room, so re-fit the row
# heights to match the new content.
# ---------------------------------------------------------------------------

$ws.Rows.Item(2).RowHeight = 348
$ws.Rows.Item(3).RowHeight = 261
$ws.Rows.Item(4).RowHeight = 409.5
$ws.Rows.Item(5).RowHeight = 409.5

# ---------------------------------------------------------------------------
# Update the saved view/selection state.
# ---------------------------------------------------------------------------

$ws.Activate() | Out-Null
$ws.Range("B3").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1

Write-Host "Updated CasesTab/SamplesTab/FilesTab/StudyFilesTab queries"
